$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete row 2 (the 2007/2008 entry with no C or E value), shifting all rows up by one
$ws.Rows("2:2").Delete()

# Step 2: clear cells that should become empty after the recompute
$ws.Cells.Item(2, 5).ClearContents()
$ws.Cells.Item(3, 5).ClearContents()
$ws.Cells.Item(4, 3).ClearContents()
$ws.Cells.Item(4, 5).ClearContents()
$ws.Cells.Item(5, 5).ClearContents()
$ws.Cells.Item(6, 3).ClearContents()
$ws.Cells.Item(6, 5).ClearContents()
$ws.Cells.Item(7, 5).ClearContents()
$ws.Cells.Item(8, 3).ClearContents()
$ws.Cells.Item(8, 5).ClearContents()
$ws.Cells.Item(9, 5).ClearContents()
$ws.Cells.Item(10, 3).ClearContents()
$ws.Cells.Item(10, 5).ClearContents()

# Step 3: set updated/recomputed values for C (y_0_forecast) and E (y_1_forecast) columns
$ws.Cells.Item(5, 3).Value = 1.834695583582535
$ws.Cells.Item(7, 3).Value = 1.767835936772144
$ws.Cells.Item(9, 3).Value = 1.074400434091038
$ws.Cells.Item(11, 5).Value = 1.274704633957136
$ws.Cells.Item(12, 5).Value = 1.37755776875883
$ws.Cells.Item(12, 3).Value = 1.097054137926201
$ws.Cells.Item(13, 5).Value = 1.404348988410131
$ws.Cells.Item(14, 5).Value = 1.329814931661888
$ws.Cells.Item(14, 3).Value = 1.385527545913412
$ws.Cells.Item(15, 5).Value = 1.269653854937691
$ws.Cells.Item(16, 5).Value = 1.355477993452414
$ws.Cells.Item(16, 3).Value = 1.296301936385214
$ws.Cells.Item(17, 5).Value = 1.253742200752095
$ws.Cells.Item(18, 5).Value = 1.386547975635688
$ws.Cells.Item(18, 3).Value = 2.441628883342295
$ws.Cells.Item(19, 5).Value = 1.805615391969595
$ws.Cells.Item(20, 5).Value = 1.833587970352424
$ws.Cells.Item(20, 3).Value = 2.565764046666463
$ws.Cells.Item(21, 5).Value = 1.661541796722577
$ws.Cells.Item(22, 5).Value = 1.553425185274571
$ws.Cells.Item(22, 3).Value = 1.376993627314671
$ws.Cells.Item(23, 5).Value = 1.485511920344451
$ws.Cells.Item(23, 3).Value = 1.263447557103259
$ws.Cells.Item(24, 5).Value = 1.635045928803081
$ws.Cells.Item(24, 3).Value = 1.362852986880547
$ws.Cells.Item(25, 5).Value = 1.815016201748643
$ws.Cells.Item(26, 5).Value = 1.600603483732033
$ws.Cells.Item(26, 3).Value = 1.868356483387124
$ws.Cells.Item(27, 5).Value = 1.745834498329324
$ws.Cells.Item(27, 3).Value = 2.117022522597423
$ws.Cells.Item(28, 5).Value = 1.944926006147352
$ws.Cells.Item(28, 3).Value = 2.244229492550187
$ws.Cells.Item(29, 5).Value = 1.810449264563152
$ws.Cells.Item(30, 5).Value = 1.679039075934385
$ws.Cells.Item(30, 3).Value = 2.005372766276703
$ws.Cells.Item(31, 5).Value = 1.76475225558832
$ws.Cells.Item(31, 3).Value = 2.149400276001101
$ws.Cells.Item(32, 5).Value = 1.674992401025355
$ws.Cells.Item(33, 5).Value = 2.128328071999674
$ws.Cells.Item(34, 5).Value = 1.827015901454287
$ws.Cells.Item(34, 3).Value = 2.273520076663971
$ws.Cells.Item(35, 5).Value = 2.131436976903012
$ws.Cells.Item(35, 3).Value = 2.453568910971748
$ws.Cells.Item(36, 5).Value = 0.5797569954799853
$ws.Cells.Item(37, 5).Value = 1.11435041103376
$ws.Cells.Item(38, 5).Value = 1.421947874507667
$ws.Cells.Item(38, 3).Value = 0.7380952286421882
$ws.Cells.Item(39, 5).Value = 1.556352278772266
$ws.Cells.Item(39, 3).Value = 0.812682184439506
$ws.Cells.Item(40, 5).Value = 1.891749670939347
$ws.Cells.Item(41, 5).Value = 0.3338851812143995
$ws.Cells.Item(42, 5).Value = 1.173782198617435
$ws.Cells.Item(42, 3).Value = -0.4131792716363547
$ws.Cells.Item(43, 5).Value = 1.820779918499094
$ws.Cells.Item(43, 3).Value = 0.9940067218177528
$ws.Cells.Item(44, 5).Value = 1.655852464312013
$ws.Cells.Item(45, 5).Value = 2.228542839642689
$ws.Cells.Item(46, 5).Value = 1.632015075917925
$ws.Cells.Item(46, 3).Value = 1.805571054927801
$ws.Cells.Item(47, 5).Value = 1.554016159863814
$ws.Cells.Item(47, 3).Value = 1.634555928116921
$ws.Cells.Item(48, 5).Value = 1.310740881193517
$ws.Cells.Item(49, 5).Value = 0.9823016603409229
$ws.Cells.Item(50, 5).Value = 1.590970324046337
$ws.Cells.Item(50, 3).Value = 0.6395223689078522
$ws.Cells.Item(51, 5).Value = 1.581524829939718
$ws.Cells.Item(51, 3).Value = 0.6231570351797
$ws.Cells.Item(52, 5).Value = 1.539180932382078

Write-Host "Done applying naive forecaster bugfix"
